# Update the "timestamp" column (Z) values for every data row (2-48) with
# the new run's timestamps, grouped by the batches that share an identical
# timestamp (mirrors how the original logging loop stamped groups of rows
# written in the same processing batch).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Z2:Z15").Value  = "2025-10-17T07:09:29.479980"
$ws.Range("Z16:Z21").Value = "2025-10-17T07:09:29.567060"
$ws.Range("Z22:Z25").Value = "2025-10-17T07:09:29.568058"
$ws.Range("Z26:Z28").Value = "2025-10-17T07:09:29.651676"
$ws.Range("Z29:Z34").Value = "2025-10-17T07:09:29.652676"
$ws.Range("Z35:Z45").Value = "2025-10-17T07:09:29.653690"
$ws.Range("Z46:Z48").Value = "2025-10-17T07:09:29.654681"
